$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Rebuild the worker mora table (rows 16-48) in reverse chronological period order,
# per "Elimina EC anteriores y se agregan nuevos, se modifica base de datos".
$ws.Range("C16").Value = "45692521"
$ws.Range("D16").Value = "MILENA DEL SOCORRO JIMENEZ GUTIERREZ"
$ws.Range("E16").Value = "2403"
$ws.Range("F16").Value = 25333
$ws.Range("G16").Value = 1048000
$ws.Range("C17").Value = "45692521"
$ws.Range("D17").Value = "MILENA DEL SOCORRO JIMENEZ GUTIERREZ"
$ws.Range("E17").Value = "2402"
$ws.Range("F17").Value = 40000
$ws.Range("G17").Value = 1048000
$ws.Range("C18").Value = "45692521"
$ws.Range("D18").Value = "MILENA DEL SOCORRO JIMENEZ GUTIERREZ"
$ws.Range("E18").Value = "2401"
$ws.Range("F18").Value = 40000
$ws.Range("G18").Value = 1048000
$ws.Range("C19").Value = "45692521"
$ws.Range("D19").Value = "MILENA DEL SOCORRO JIMENEZ GUTIERREZ"
$ws.Range("E19").Value = "2312"
$ws.Range("F19").Value = 40000
$ws.Range("G19").Value = 1048000
$ws.Range("C20").Value = "45692521"
$ws.Range("D20").Value = "MILENA DEL SOCORRO JIMENEZ GUTIERREZ"
$ws.Range("E20").Value = "2311"
$ws.Range("F20").Value = 40000
$ws.Range("G20").Value = 1048000
$ws.Range("C21").Value = "45692521"
$ws.Range("D21").Value = "MILENA DEL SOCORRO JIMENEZ GUTIERREZ"
$ws.Range("E21").Value = "2310"
$ws.Range("F21").Value = 40000
$ws.Range("G21").Value = 1048000
$ws.Range("C22").Value = "45692521"
$ws.Range("D22").Value = "MILENA DEL SOCORRO JIMENEZ GUTIERREZ"
$ws.Range("E22").Value = "2309"
$ws.Range("F22").Value = 40000
$ws.Range("G22").Value = 1048000
$ws.Range("C23").Value = "45692521"
$ws.Range("D23").Value = "MILENA DEL SOCORRO JIMENEZ GUTIERREZ"
$ws.Range("E23").Value = "2308"
$ws.Range("F23").Value = 40000
$ws.Range("G23").Value = 1048000
$ws.Range("C24").Value = "45692521"
$ws.Range("D24").Value = "MILENA DEL SOCORRO JIMENEZ GUTIERREZ"
$ws.Range("E24").Value = "2307"
$ws.Range("F24").Value = 40000
$ws.Range("G24").Value = 1048000
$ws.Range("C25").Value = "45692521"
$ws.Range("D25").Value = "MILENA DEL SOCORRO JIMENEZ GUTIERREZ"
$ws.Range("E25").Value = "2306"
$ws.Range("F25").Value = 40000
$ws.Range("G25").Value = 1048000
$ws.Range("C26").Value = "45692521"
$ws.Range("D26").Value = "MILENA DEL SOCORRO JIMENEZ GUTIERREZ"
$ws.Range("E26").Value = "2305"
$ws.Range("F26").Value = 40000
$ws.Range("G26").Value = 1048000
$ws.Range("C27").Value = "45692521"
$ws.Range("D27").Value = "MILENA DEL SOCORRO JIMENEZ GUTIERREZ"
$ws.Range("E27").Value = "2304"
$ws.Range("F27").Value = 40000
$ws.Range("G27").Value = 1048000
$ws.Range("C28").Value = "45692521"
$ws.Range("D28").Value = "MILENA DEL SOCORRO JIMENEZ GUTIERREZ"
$ws.Range("E28").Value = "2303"
$ws.Range("F28").Value = 40000
$ws.Range("G28").Value = 1048000
$ws.Range("C29").Value = "45692521"
$ws.Range("D29").Value = "MILENA DEL SOCORRO JIMENEZ GUTIERREZ"
$ws.Range("E29").Value = "2302"
$ws.Range("F29").Value = 40000
$ws.Range("G29").Value = 1048000
$ws.Range("C30").Value = "45692521"
$ws.Range("D30").Value = "MILENA DEL SOCORRO JIMENEZ GUTIERREZ"
$ws.Range("E30").Value = "2301"
$ws.Range("F30").Value = 40000
$ws.Range("G30").Value = 1048000
$ws.Range("C31").Value = "45692521"
$ws.Range("D31").Value = "MILENA DEL SOCORRO JIMENEZ GUTIERREZ"
$ws.Range("E31").Value = "2212"
$ws.Range("F31").Value = 40000
$ws.Range("G31").Value = 1048000
$ws.Range("C32").Value = "45692521"
$ws.Range("D32").Value = "MILENA DEL SOCORRO JIMENEZ GUTIERREZ"
$ws.Range("E32").Value = "2211"
$ws.Range("F32").Value = 40000
$ws.Range("G32").Value = 1048000
$ws.Range("C33").Value = "45692521"
$ws.Range("D33").Value = "MILENA DEL SOCORRO JIMENEZ GUTIERREZ"
$ws.Range("E33").Value = "2210"
$ws.Range("F33").Value = 40000
$ws.Range("G33").Value = 1048000
$ws.Range("C34").Value = "45692521"
$ws.Range("D34").Value = "MILENA DEL SOCORRO JIMENEZ GUTIERREZ"
$ws.Range("E34").Value = "2209"
$ws.Range("F34").Value = 40000
$ws.Range("G34").Value = 1048000
$ws.Range("C35").Value = "45692521"
$ws.Range("D35").Value = "MILENA DEL SOCORRO JIMENEZ GUTIERREZ"
$ws.Range("E35").Value = "2208"
$ws.Range("F35").Value = 40000
$ws.Range("G35").Value = 1048000
$ws.Range("C36").Value = "45692521"
$ws.Range("D36").Value = "MILENA DEL SOCORRO JIMENEZ GUTIERREZ"
$ws.Range("E36").Value = "2207"
$ws.Range("F36").Value = 40000
$ws.Range("G36").Value = 1048000
$ws.Range("C37").Value = "1051890322"
$ws.Range("D37").Value = "ROSMERIS HERRERA GONZALEZ"
$ws.Range("E37").Value = "2403"
$ws.Range("F37").Value = 29387
$ws.Range("G37").Value = 1160000
$ws.Range("C38").Value = "1051890322"
$ws.Range("D38").Value = "ROSMERIS HERRERA GONZALEZ"
$ws.Range("E38").Value = "2402"
$ws.Range("F38").Value = 46400
$ws.Range("G38").Value = 1160000
$ws.Range("C39").Value = "1051890322"
$ws.Range("D39").Value = "ROSMERIS HERRERA GONZALEZ"
$ws.Range("E39").Value = "2401"
$ws.Range("F39").Value = 46400
$ws.Range("G39").Value = 1160000
$ws.Range("C40").Value = "1051890322"
$ws.Range("D40").Value = "ROSMERIS HERRERA GONZALEZ"
$ws.Range("E40").Value = "2312"
$ws.Range("F40").Value = 46400
$ws.Range("G40").Value = 1160000
$ws.Range("C41").Value = "1051890322"
$ws.Range("D41").Value = "ROSMERIS HERRERA GONZALEZ"
$ws.Range("E41").Value = "2311"
$ws.Range("F41").Value = 46400
$ws.Range("G41").Value = 1160000
$ws.Range("C42").Value = "1051890322"
$ws.Range("D42").Value = "ROSMERIS HERRERA GONZALEZ"
$ws.Range("E42").Value = "2310"
$ws.Range("F42").Value = 46400
$ws.Range("G42").Value = 1160000
$ws.Range("C43").Value = "1051890322"
$ws.Range("D43").Value = "ROSMERIS HERRERA GONZALEZ"
$ws.Range("E43").Value = "2309"
$ws.Range("F43").Value = 46400
$ws.Range("G43").Value = 1160000
$ws.Range("C44").Value = "1051890322"
$ws.Range("D44").Value = "ROSMERIS HERRERA GONZALEZ"
$ws.Range("E44").Value = "2308"
$ws.Range("F44").Value = 46400
$ws.Range("G44").Value = 1160000
$ws.Range("C45").Value = "1051890322"
$ws.Range("D45").Value = "ROSMERIS HERRERA GONZALEZ"
$ws.Range("E45").Value = "2307"
$ws.Range("F45").Value = 46400
$ws.Range("G45").Value = 1160000
$ws.Range("C46").Value = "1051890322"
$ws.Range("D46").Value = "ROSMERIS HERRERA GONZALEZ"
$ws.Range("E46").Value = "2306"
$ws.Range("F46").Value = 46400
$ws.Range("G46").Value = 1160000
$ws.Range("C47").Value = "1051890322"
$ws.Range("D47").Value = "ROSMERIS HERRERA GONZALEZ"
$ws.Range("E47").Value = "2305"
$ws.Range("F47").Value = 46400
$ws.Range("G47").Value = 1160000
$ws.Range("C48").Value = "1051890322"
$ws.Range("D48").Value = "ROSMERIS HERRERA GONZALEZ"
$ws.Range("E48").Value = "2304"
$ws.Range("F48").Value = 46400
$ws.Range("G48").Value = 1160000
